# Rename the *img sheets to img* and make the last one (imge, formerly eimg)
# the active/selected sheet, matching the author's commit "Change names from *img to img*".

$wb = $excel.ActiveWorkbook

$renames = @{
    "himg" = "imgh"
    "timg" = "imgt"
    "simg" = "imgs"
    "gimg" = "imgg"
    "wimg" = "imgw"
    "bimg" = "imgb"
    "eimg" = "imge"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# Activate the last renamed sheet (now named "imge") so it becomes the
# workbook's active tab and its sheetView gets tabSelected="true".
$wb.Worksheets.Item("imge").Activate()
